$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "245.64"
Set-TextValue $ws.Range("D3") "22.05"
Set-TextValue $ws.Range("D4") "5.406"
Set-TextValue $ws.Range("D5") "0.05855"
Set-TextValue $ws.Range("D8") "0.8185"
Set-TextValue $ws.Range("D9") "1.018"
Set-TextValue $ws.Range("D11") "0.07440"
Set-TextValue $ws.Range("D12") "0.03446"
Set-TextValue $ws.Range("D13") "0.03044"
Set-TextValue $ws.Range("D14") "4.184"
Set-TextValue $ws.Range("D15") "0.09394"
Set-TextValue $ws.Range("D16") "0.001602"
Set-TextValue $ws.Range("D17") "0.04836"
Set-TextValue $ws.Range("D18") "0.0005895"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue $ws.Range("D19") "0.006024"
Set-TextValue $ws.Range("D20") "0.004100"
Set-TextValue $ws.Range("D21") "0.0009976"
Set-TextValue $ws.Range("D22") "0.0001501"
Set-TextValue $ws.Range("D23") "3.696"
Set-TextValue $ws.Range("D24") "2.215"
Set-TextValue $ws.Range("D25") "0.3237"
$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"
Set-TextValue $ws.Range("D40") "0.03864"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006466"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.002602"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws.Range("D44") "0.006254"
Set-TextValue $ws.Range("D45") "0.00005627"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("D47") "0.4204"
Set-TextValue $ws.Range("D48") "0.1421"
Set-TextValue $ws.Range("D49") "0.00002102"
Set-TextValue $ws.Range("D50") "0.01011"
